$d = $word.ActiveDocument

$d.Content.Find.Execute("909÷6=151, 3", $true, $false, $false, $false, $false, $true, 1, $false, "458÷2=229, 0", 2) | Out-Null
$d.Content.Find.Execute("220÷3=73, 1", $true, $false, $false, $false, $false, $true, 1, $false, "511÷5=102, 1", 2) | Out-Null
$d.Content.Find.Execute("568÷6=94, 4", $true, $false, $false, $false, $false, $true, 1, $false, "150÷3=50, 0", 2) | Out-Null
$d.Content.Find.Execute("276÷2=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "941÷3=313, 2", 2) | Out-Null
$d.Content.Find.Execute("234÷6=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "920÷5=184, 0", 2) | Out-Null
$d.Content.Find.Execute("288÷8=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "217÷6=36, 1", 2) | Out-Null
$d.Content.Find.Execute("882÷5=176, 2", $true, $false, $false, $false, $false, $true, 1, $false, "967÷4=241, 3", 2) | Out-Null
$d.Content.Find.Execute("525÷5=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "384÷7=54, 6", 2) | Out-Null
$d.Content.Find.Execute("121÷2=60, 1", $true, $false, $false, $false, $false, $true, 1, $false, "597÷6=99, 3", 2) | Out-Null
$d.Content.Find.Execute("647÷8=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "306÷6=51, 0", 2) | Out-Null
$d.Content.Find.Execute("308÷5=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "562÷5=112, 2", 2) | Out-Null
$d.Content.Find.Execute("565÷2=282, 1", $true, $false, $false, $false, $false, $true, 1, $false, "494÷7=70, 4", 2) | Out-Null
$d.Content.Find.Execute("780÷2=390, 0", $true, $false, $false, $false, $false, $true, 1, $false, "459÷6=76, 3", 2) | Out-Null
$d.Content.Find.Execute("926÷4=231, 2", $true, $false, $false, $false, $false, $true, 1, $false, "265÷9=29, 4", 2) | Out-Null
$d.Content.Find.Execute("515÷3=171, 2", $true, $false, $false, $false, $false, $true, 1, $false, "794÷7=113, 3", 2) | Out-Null
$d.Content.Find.Execute("935÷3=311, 2", $true, $false, $false, $false, $false, $true, 1, $false, "363÷5=72, 3", 2) | Out-Null
$d.Content.Find.Execute("809÷4=202, 1", $true, $false, $false, $false, $false, $true, 1, $false, "657÷3=219, 0", 2) | Out-Null
$d.Content.Find.Execute("278÷3=92, 2", $true, $false, $false, $false, $false, $true, 1, $false, "985÷3=328, 1", 2) | Out-Null
$d.Content.Find.Execute("326÷8=40, 6", $true, $false, $false, $false, $false, $true, 1, $false, "766÷3=255, 1", 2) | Out-Null
$d.Content.Find.Execute("751÷7=107, 2", $true, $false, $false, $false, $false, $true, 1, $false, "407÷9=45, 2", 2) | Out-Null
$d.Content.Find.Execute("468÷6=78, 0", $true, $false, $false, $false, $false, $true, 1, $false, "206÷8=25, 6", 2) | Out-Null
$d.Content.Find.Execute("961÷7=137, 2", $true, $false, $false, $false, $false, $true, 1, $false, "151÷2=75, 1", 2) | Out-Null
$d.Content.Find.Execute("768÷4=192, 0", $true, $false, $false, $false, $false, $true, 1, $false, "185÷2=92, 1", 2) | Out-Null
$d.Content.Find.Execute("915÷9=101, 6", $true, $false, $false, $false, $false, $true, 1, $false, "369÷7=52, 5", 2) | Out-Null
$d.Content.Find.Execute("233÷4=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "210÷3=70, 0", 2) | Out-Null
